$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.633.65"
$ws.Range("E2").Value = "  +0.64%  "

# Row 3
$ws.Range("D3").Value = "1.802.03"
$ws.Range("E3").Value = "  -0.96%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.56%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6
$ws.Range("E6").Value = "  +0.50%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5318"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3771"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.24%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07492"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.04%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.65%  "

# Row 12
$ws.Range("E12").Value = "  +0.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.148"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.333"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

# Row 16
$ws.Range("D16").Value = "1.805.15"
$ws.Range("E16").Value = "  -0.38%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.99%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06452"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.00%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.905"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.81%  "

# Row 23
$ws.Range("D23").Value = "28.667.24"
$ws.Range("E23").Value = "  +0.72%  "

# Row 24
$ws.Range("E24").Value = "  -2.79%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.094"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.85%  "

# Row 28
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.012.09"
$ws.Range("E28").Value = "  -0.51%  "

# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.355"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.096"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.15%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.89%  "

# Row 33
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.694"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.60%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.631"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2248"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.61%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06404"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.88%  "

# Row 37
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.802"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02299"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.05%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.022"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.229"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.34%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.31%  "

# Row 43
$ws.Range("E43").Value = "  +0.50%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.405"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.24%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.690"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5845"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.938"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.146"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06887"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
